$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '69.308.47'
$ws.Range('E2').Value = '  +1.54%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.898.35'
$ws.Range('E3').Value = '  -0.48%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '526.15'
$ws.Range('E5').Value = '  +8.49%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '143.03'
$ws.Range('E6').Value = '  -2.22%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.610'
$ws.Range('E7').Value = '  -1.91%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.718'
$ws.Range('E9').Value = '  -2.65%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.171'
$ws.Range('E10').Value = '  +1.22%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0000331'
$ws.Range('E11').Value = '  -4.38%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '41.91'
$ws.Range('E12').Value = '  -2.99%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.508.25'
$ws.Range('E13').Value = '  -0.68%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '10.18'
$ws.Range('E14').Value = '  -5.00%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.896.41'
$ws.Range('E15').Value = '  -0.53%  '
$ws.Range('E16').Value = '  +7.54%  '
$ws.Range('E17').Value = '  -0.65%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '13.81'
$ws.Range('E18').Value = '  -2.99%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '19.67'
$ws.Range('E19').Value = '  -2.93%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '69.211.81'
$ws.Range('E20').Value = '  +1.36%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '424.78'
$ws.Range('E21').Value = '  -1.48%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '3.34'
$ws.Range('E22').Value = '  -5.29%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '14.13'
$ws.Range('E23').Value = '  -6.36%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '87.91'
$ws.Range('E24').Value = '  -1.49%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '4.03'
$ws.Range('E25').Value = '  +8.56%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '11.64'
$ws.Range('E26').Value = '  -0.95%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.53'
$ws.Range('E27').Value = '  -5.83%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '36.12'
$ws.Range('E28').Value = '  -4.42%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '695.15'
$ws.Range('E29').Value = '  -3.35%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '13.10'
$ws.Range('E30').Value = '  -4.83%  '
$ws.Range('E31').Value = '  -4.04%  '
$ws.Range('E32').Value = '  -4.37%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '68.08'
$ws.Range('E33').Value = '  +11.37%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.441'
$ws.Range('E34').Value = '  +9.39%  '
$ws.Range('E35').Value = '  -5.03%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '40.00'
$ws.Range('E36').Value = '  -4.07%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.0₃0833'
$ws.Range('E37').Value = '  -6.35%  '
$ws.Range('E38').Value = '  +2.82%  '
$ws.Range('E39').Value = '  +0.31%  '
$ws.Range('E40').Value = '  -0.14%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0479'
$ws.Range('E41').Value = '  -2.16%  '
$ws.Range('E42').Value = '  -10.46%  '
$ws.Range('E43').Value = '  +0.30%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.96'
$ws.Range('E44').Value = '  -5.44%  '
$ws.Range('E45').Value = '  -1.26%  '
$ws.Range('E46').Value = '  -1.95%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.03'
$ws.Range('E47').Value = '  +8.03%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '26.70'
$ws.Range('E48').Value = '  +5.49%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '3.27'
$ws.Range('E49').Value = '  -4.81%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0₆0340'
$ws.Range('E50').Value = '  +1.56%  '
$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '142.52'
$ws.Range('E51').Value = '  -2.01%  '
